$d = $word.ActiveDocument

# 1. Update activation date
$null = $d.Content.Find.Execute(
    "Ativação: 01/01/2022", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2025", 2)

# 2. Add trailing period to the English "Objetivos" paragraph
$oldObjEn = "Understand the theoretical bases of qualitative analytical chemistry of environmental interest"
$newObjEn = "Understand the theoretical bases of qualitative analytical chemistry of environmental interest."
$null = $d.Content.Find.Execute($oldObjEn, $true, $false, $false, $false, $false,
    $true, 1, $false, $newObjEn, 2)

# 3. Insert four new "Docente(s) Responsável(eis)" list entries before the existing one
$rng = $d.Content
$found = $rng.Find.Execute("7455355 - Robson da Silva Rocha", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if ($found) {
    $insertPoint = $d.Range($rng.Start, $rng.Start)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
        '<w:r><w:t>7043088 - Ana Karine Furtado de Carvalho</w:t><w:br/></w:r>' +
        '<w:r><w:t>7926291 - Célia Regina Tomachuk dos Santos Catuogno</w:t><w:br/></w:r>' +
        '<w:r><w:t>4893449 - Débora Souza Alvim</w:t><w:br/></w:r>' +
        '<w:r><w:t>8855158 - Morun Bernardino Neto</w:t><w:br/></w:r>' +
        '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertPoint.InsertXML($xml)
}

# 4. Shorten "Programa resumido" paragraph (Portuguese)
$oldResumoPt = "Introdução à análise qualitativa, indicando suas aplicabilidades e limitações. Uso das técnicas qualitativas para análise dos principais íons de importância ambiental. Análise de sólidos, partículas, sedimentos. Estudos de amostras de importância ambiental."
$newResumoPt = "Introdução à análise qualitativa, indicando suas aplicabilidades e limitações. Uso das técnicas qualitativas para análise dos principais íons de importância ambiental."
$null = $d.Content.Find.Execute($oldResumoPt, $true, $false, $false, $false, $false,
    $true, 1, $false, $newResumoPt, 2)

# 5. Shorten "Programa resumido" paragraph (English)
$oldResumoEn = "Introduction to qualitative analysis, indicating its applicability and limitations. Use of qualitative techniques to analyze the main ions of environmental importance. Analysis of solids, particles, sediments. Studies of samples of environmental importance."
$newResumoEn = "Introduction to qualitative analysis, indicating its applicability and limitations. Use of qualitative techniques to analyze the main ions of environmental importance."
$null = $d.Content.Find.Execute($oldResumoEn, $true, $false, $false, $false, $false,
    $true, 1, $false, $newResumoEn, 2)

# 6. Shorten "Programa" paragraph (Portuguese)
$oldProgPt = "- Revisão das regras de segurança laboratorial - Introdução à análise qualitativa: Definições, objetivos e limitações. - Análise de sólidos, partículas, sedimentos.- Identificação dos cátions do grupo I (K+, Na+ e NH4+); grupo II (Mg2+, Ca2+ e Ba2+); grupo III (Al3+, Fe3+, Mn2+).- Estudo dos ânions e suas aplicações em análises ambientais (Cl e suas espécies, SO42-, CO32-, S2-, NO3-).- Análise gravimétrica: fundamentos e formação de precipitados.- Análises dos principais cátions e ânions em amostras conhecidas e desconhecidas para os alunos- Análise de metais em solo, água ou outras amostras ambientais importantes"
$newProgPt = "- Revisão das regras de segurança laboratorial - Introdução à análise qualitativa: Definições, objetivos e limitações. - Identificação dos cátions do grupo I (K+, Na+ e NH4+); grupo II (Mg2+, Ca2+ e Ba2+); grupo III (Al3+, Fe3+, Mn2+).- Estudo dos ânions e suas aplicações em análises ambientais (Cl e suas espécies, SO42-, CO32-, S2-, NO3-)."
$null = $d.Content.Find.Execute($oldProgPt, $true, $false, $false, $false, $false,
    $true, 1, $false, $newProgPt, 2)

# 7. Shorten "Programa" paragraph (English)
$oldProgEn = "- Review of laboratory safety rules- Introduction to qualitative analysis: Definitions, objectives and limitations.- Analysis of solids, particles, sediments.- Identification of group I cations (K+, Na+ and NH4+); group II (Mg2+, Ca2+ and Ba2+); group III (Al3+, Fe3+, Fe2+, Mn2+).- Study of anions and their applications in environmental analysis (Cl and its species, SO42-, CO32-, S2-, NO3-).- Gravimetric analysis: fundamentals and precipitate formation.- Analysis of the main cations and anions in known and unknown samples for students- Analysis of metals in soil, water or other important environmental samples"
$newProgEn = "- Review of laboratory safety rules- Introduction to qualitative analysis: Definitions, objectives and limitations.- Identification of group I cations (K+, Na+ and NH4+); group II (Mg2+, Ca2+ and Ba2+); group III (Al3+, Fe3+, Fe2+, Mn2+).- Study of anions and their applications in environmental analysis (Cl and its species, SO42-, CO32-, S2-, NO3-)."
$null = $d.Content.Find.Execute($oldProgEn, $true, $false, $false, $false, $false,
    $true, 1, $false, $newProgEn, 2)

Write-Host "Done"
